$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Tabelle1")

# The "Datum" column (A) held full dates like "31.12.2010"; strip the
# "31.12." day/month prefix so only the year remains, e.g. "2010".
for ($r = 2; $r -le 15; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    $cell.Value = $old -replace '^31\.12\.', ''
}

# Update the view: scrolled down a bit with a different active cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("D13").Select()
